$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the pre-1949 / pre-1969 "RES1/RES2" split labels to their new
# construction-year-based names (filter_tek_list small update).
$ws.Range("A2").Value = "PRE_TEK49_RES_1950"
$ws.Range("A3").Value = "PRE_TEK49_RES_1940"
$ws.Range("A7").Value = "TEK69_RES_1976"
$ws.Range("A8").Value = "TEK69_RES_1986"

# Add a new (currently empty) percentage-formatted column next to every
# "RES"/"COM" TEK-period row (rows 3-11), matching the new cellXfs entry
# (numFmtId 9, i.e. "0%").
$percentRows = 3..11
foreach ($r in $percentRows) {
    $ws.Cells.Item($r, 3).NumberFormat = "0%"
}

# Column A grew wider text ("PRE_TEK49_RES_1950" / "TEK69_RES_1976", etc.)
# and column C needs its own best-fit width; reflect both explicitly
# (values tuned so the exported character-width lands as close as this
# interop layer's column-width quantization allows to Excel's own
# best-fit result of ~19.71 / ~10.43).
$ws.Columns.Item(1).ColumnWidth = 18.76
$ws.Columns.Item(3).ColumnWidth = 9.6

# Selection moved to A3 after the edit.
$ws.Range("A3").Select() | Out-Null
